# Updates crypto price/volume data (Price + Volume(1h) columns) to refresh
# the snapshot from the latest run of the cryptos-list GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds plain-text numbers (e.g. "512.00", "3.028.47")
# -- force Text format first so Excel does not silently reinterpret them
# as real numbers (which would drop trailing zeros / thousands dots).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '56.891.77'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '3.040.48'
$ws.Range('E3').Value = '  +2.70%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '512.00'
$ws.Range('E5').Value = '  +3.06%  '
$ws.Range('D6').Value = '140.96'
$ws.Range('E6').Value = '  +4.74%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('D8').Value = '0.432'
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('D9').Value = '7.16'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').Value = '0.109'
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('D11').Value = '0.371'
$ws.Range('E11').Value = '  +5.47%  '
$ws.Range('D12').Value = '3.549.38'
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '25.44'
$ws.Range('E14').Value = '  -1.94%  '
$ws.Range('D15').Value = '0.0000164'
$ws.Range('E15').Value = '  +4.37%  '
$ws.Range('D16').Value = '56.867.18'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').Value = '3.033.73'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D18').Value = '5.95'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').Value = '13.21'
$ws.Range('E19').Value = '  +5.88%  '
$ws.Range('D20').Value = '8.10'
$ws.Range('E20').Value = '  +4.56%  '
$ws.Range('D21').Value = '335.10'
$ws.Range('E21').Value = '  +5.54%  '
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').Value = '0.503'
$ws.Range('E23').Value = '  +3.50%  '
$ws.Range('D24').Value = '64.79'
$ws.Range('E24').Value = '  +3.38%  '
$ws.Range('D25').Value = '3.156.25'
$ws.Range('E25').Value = '  +2.49%  '
$ws.Range('D26').Value = '0.167'
$ws.Range('E26').Value = '  +2.78%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').Value = '0.0₃0941'
$ws.Range('E28').Value = '  +9.13%  '
$ws.Range('D29').Value = '6.47'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').Value = '6.81'
$ws.Range('E30').Value = '  -2.82%  '
$ws.Range('D31').Value = '1.80'
$ws.Range('E31').Value = '  +3.17%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '20.52'
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.17'
$ws.Range('E33').Value = '  +3.50%  '
$ws.Range('D34').Value = '152.96'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('D35').Value = '4.51'
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('D36').Value = '27.21'
$ws.Range('E36').Value = '  +13.66%  '
$ws.Range('D37').Value = '5.85'
$ws.Range('E37').Value = '  +2.90%  '
$ws.Range('E38').Value = '  +2.28%  '
$ws.Range('D39').Value = '0.0666'
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('D40').Value = '3.074.19'
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('D41').Value = '36.64'
$ws.Range('E41').Value = '  -1.69%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = '3.82'
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('D44').Value = '0.661'
$ws.Range('E44').Value = '  +3.29%  '
$ws.Range('D45').Value = '2.211.87'
$ws.Range('E45').Value = '  +2.73%  '
$ws.Range('D46').Value = '1.36'
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('D47').Value = '0.0245'
$ws.Range('E47').Value = '  +5.79%  '
$ws.Range('D48').Value = '0.942'
$ws.Range('E48').Value = '  +2.00%  '
$ws.Range('D49').Value = '5.87'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').Value = '19.87'
$ws.Range('E50').Value = '  +4.79%  '
$ws.Range('D51').Value = '0.0858'
$ws.Range('E51').Value = '  +0.80%  '
